# Applies the "weitere Vitalparameter und englische Version des IGs hinzugefügt" edit:
#  - Title translated to English
#  - Date updated
#  - Contact URL updated
#  - Jurisdiction value set to "Germany"
#  - Description translated to English

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B5").Value = "T-CABS ValueSet Ventilation Type"
$ws.Range("B8").Value = "2025-11-19T11:55:29+01:00"
$ws.Range("B10").Value = "BIH-CEI (https://www.bihealth.org/)"
$ws.Range("B11").Value = "Germany"
$ws.Range("B12").Value = "This ValueSet contains codes to represent the different types of ventilation"
